$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to refreshed values from the scraper run.
# For column D values that look like plain numbers, force text storage (NumberFormat "@")
# so exact formatting (e.g. trailing zeros "19.00") survives, then ClearFormats() so the
# cell keeps its original (default) style -- only the stored value changes.
$ws.Range('D2').Value = '68.885.99'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '2.441.94'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.11'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.46'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.512'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('E9').Value = '  +10.81%  '
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.332'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.59'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.34%  '
$ws.Range('E13').Value = '  +5.77%  '
$ws.Range('D14').Value = '68.775.13'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '2.893.01'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.33'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '2.445.17'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.57'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.55'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.96'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.93'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.14'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').Value = '2.570.31'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('D29').Value = '0.0₃0826'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '431.61'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('E33').Value = '  +2.76%  '
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '157.99'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.02'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('E41').Value = '  +3.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.37'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.07'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.36'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '129.74'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.556'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0925'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('E51').Value = '  +1.00%  '
